# Update "想去人数" (interested-attendee count) figures in the 展览 (rId1)
# and 全部类型 (rId4) sheets to match the newly scraped numbers.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibition) ---
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F2").Value = 332
$wsExpo.Range("F3").Value = 8530
$wsExpo.Range("F4").Value = 6238
$wsExpo.Range("F5").Value = 546
$wsExpo.Range("F9").Value = 334
$wsExpo.Range("F10").Value = 1223
$wsExpo.Range("F11").Value = 88

# --- Sheet "全部类型" (All types) ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 332
$wsAll.Range("F3").Value = 8530
$wsAll.Range("F4").Value = 6238
$wsAll.Range("F5").Value = 546
$wsAll.Range("F9").Value = 334
$wsAll.Range("F14").Value = 1223
$wsAll.Range("F15").Value = 88
